# Recompute ligand/receptor average & total expression (TPM-based), derived
# specificity scores, and edge weights for the Cadm3-Cadm3 sheet using the
# updated TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"9.426699666666666"
$ws.Range("H2").Value = [double]"28.280099"
$ws.Range("I2").Value = [double]"0.4880118193702016"
$ws.Range("J2").Value = [double]"0.4880118193702015"
$ws.Range("M2").Value = [double]"9.426699666666666"
$ws.Range("N2").Value = [double]"28.280099"
$ws.Range("O2").Value = [double]"0.4880118193702016"
$ws.Range("P2").Value = [double]"0.4880118193702015"
$ws.Range("Q2").Value = [double]"88.86266660553343"
$ws.Range("R2").Value = [double]"799.763999449801"
$ws.Range("S2").Value = [double]"0.2381555358450143"
$ws.Range("T2").Value = [double]"0.2381555358450142"

$ws.Range("G3").Value = [double]"9.426699666666666"
$ws.Range("H3").Value = [double]"28.280099"
$ws.Range("I3").Value = [double]"0.4880118193702016"
$ws.Range("J3").Value = [double]"0.4880118193702015"
$ws.Range("M3").Value = [double]"7.983522666666666"
$ws.Range("O3").Value = [double]"0.4132998355002127"
$ws.Range("P3").Value = [double]"0.4132998355002127"
$ws.Range("Q3").Value = [double]"75.25827046069243"
$ws.Range("R3").Value = [double]"677.324434146232"
$ws.Range("S3").Value = [double]"0.2016952046678638"
$ws.Range("T3").Value = [double]"0.2016952046678638"

$ws.Range("G4").Value = [double]"9.426699666666666"
$ws.Range("H4").Value = [double]"28.280099"
$ws.Range("I4").Value = [double]"0.4880118193702016"
$ws.Range("J4").Value = [double]"0.4880118193702015"
$ws.Range("M4").Value = [double]"1.816582333333334"
$ws.Range("N4").Value = [double]"5.449747"
$ws.Range("O4").Value = [double]"0.09404284435416221"
$ws.Range("P4").Value = [double]"0.0940428443541622"
$ws.Range("Q4").Value = [double]"17.12437607610589"
$ws.Range("R4").Value = [double]"154.119384684953"
$ws.Range("S4").Value = [double]"0.0458940195720234"
$ws.Range("T4").Value = [double]"0.04589401957202338"

$ws.Range("G5").Value = [double]"9.426699666666666"
$ws.Range("H5").Value = [double]"28.280099"
$ws.Range("I5").Value = [double]"0.4880118193702016"
$ws.Range("J5").Value = [double]"0.4880118193702015"
$ws.Range("M5").Value = [double]"0.08973500000000001"
$ws.Range("N5").Value = [double]"0.269205"
$ws.Range("O5").Value = [double]"0.004645500775423563"
$ws.Range("P5").Value = [double]"0.004645500775423562"
$ws.Range("Q5").Value = [double]"0.8459048945883334"
$ws.Range("R5").Value = [double]"7.613144051295"
$ws.Range("S5").Value = [double]"0.002267059285300136"
$ws.Range("T5").Value = [double]"0.002267059285300135"

$ws.Range("G6").Value = [double]"7.983522666666666"
$ws.Range("I6").Value = [double]"0.4132998355002127"
$ws.Range("J6").Value = [double]"0.4132998355002127"
$ws.Range("M6").Value = [double]"9.426699666666666"
$ws.Range("N6").Value = [double]"28.280099"
$ws.Range("O6").Value = [double]"0.4880118193702016"
$ws.Range("P6").Value = [double]"0.4880118193702015"
$ws.Range("Q6").Value = [double]"75.25827046069243"
$ws.Range("R6").Value = [double]"677.324434146232"
$ws.Range("S6").Value = [double]"0.2016952046678638"
$ws.Range("T6").Value = [double]"0.2016952046678638"

$ws.Range("G7").Value = [double]"7.983522666666666"
$ws.Range("I7").Value = [double]"0.4132998355002127"
$ws.Range("J7").Value = [double]"0.4132998355002127"
$ws.Range("M7").Value = [double]"7.983522666666666"
$ws.Range("O7").Value = [double]"0.4132998355002127"
$ws.Range("P7").Value = [double]"0.4132998355002127"
$ws.Range("Q7").Value = [double]"63.73663416918043"
$ws.Range("R7").Value = [double]"573.6297075226239"
$ws.Range("S7").Value = [double]"0.1708167540245029"
$ws.Range("T7").Value = [double]"0.1708167540245028"

$ws.Range("G8").Value = [double]"7.983522666666666"
$ws.Range("I8").Value = [double]"0.4132998355002127"
$ws.Range("J8").Value = [double]"0.4132998355002127"
$ws.Range("M8").Value = [double]"1.816582333333334"
$ws.Range("N8").Value = [double]"5.449747"
$ws.Range("O8").Value = [double]"0.09404284435416221"
$ws.Range("P8").Value = [double]"0.0940428443541622"
$ws.Range("Q8").Value = [double]"14.50272623403289"
$ws.Range("R8").Value = [double]"130.524536106296"
$ws.Range("S8").Value = [double]"0.03886789210154735"
$ws.Range("T8").Value = [double]"0.03886789210154734"

$ws.Range("G9").Value = [double]"7.983522666666666"
$ws.Range("I9").Value = [double]"0.4132998355002127"
$ws.Range("J9").Value = [double]"0.4132998355002127"
$ws.Range("M9").Value = [double]"0.08973500000000001"
$ws.Range("N9").Value = [double]"0.269205"
$ws.Range("O9").Value = [double]"0.004645500775423563"
$ws.Range("P9").Value = [double]"0.004645500775423562"
$ws.Range("Q9").Value = [double]"0.7164014064933333"
$ws.Range("R9").Value = [double]"6.44761265844"
$ws.Range("S9").Value = [double]"0.001919984706298669"
$ws.Range("T9").Value = [double]"0.001919984706298669"

$ws.Range("G10").Value = [double]"1.816582333333334"
$ws.Range("H10").Value = [double]"5.449747"
$ws.Range("I10").Value = [double]"0.09404284435416221"
$ws.Range("J10").Value = [double]"0.0940428443541622"
$ws.Range("M10").Value = [double]"9.426699666666666"
$ws.Range("N10").Value = [double]"28.280099"
$ws.Range("O10").Value = [double]"0.4880118193702016"
$ws.Range("P10").Value = [double]"0.4880118193702015"
$ws.Range("Q10").Value = [double]"17.12437607610589"
$ws.Range("R10").Value = [double]"154.119384684953"
$ws.Range("S10").Value = [double]"0.0458940195720234"
$ws.Range("T10").Value = [double]"0.04589401957202338"

$ws.Range("G11").Value = [double]"1.816582333333334"
$ws.Range("H11").Value = [double]"5.449747"
$ws.Range("I11").Value = [double]"0.09404284435416221"
$ws.Range("J11").Value = [double]"0.0940428443541622"
$ws.Range("M11").Value = [double]"7.983522666666666"
$ws.Range("O11").Value = [double]"0.4132998355002127"
$ws.Range("P11").Value = [double]"0.4132998355002127"
$ws.Range("Q11").Value = [double]"14.50272623403289"
$ws.Range("R11").Value = [double]"130.524536106296"
$ws.Range("S11").Value = [double]"0.03886789210154735"
$ws.Range("T11").Value = [double]"0.03886789210154734"

$ws.Range("G12").Value = [double]"1.816582333333334"
$ws.Range("H12").Value = [double]"5.449747"
$ws.Range("I12").Value = [double]"0.09404284435416221"
$ws.Range("J12").Value = [double]"0.0940428443541622"
$ws.Range("M12").Value = [double]"1.816582333333334"
$ws.Range("N12").Value = [double]"5.449747"
$ws.Range("O12").Value = [double]"0.09404284435416221"
$ws.Range("P12").Value = [double]"0.0940428443541622"
$ws.Range("Q12").Value = [double]"3.299971373778778"
$ws.Range("R12").Value = [double]"29.699742364009"
$ws.Range("S12").Value = [double]"0.00884405657422118"
$ws.Range("T12").Value = [double]"0.008844056574221177"

$ws.Range("G13").Value = [double]"1.816582333333334"
$ws.Range("H13").Value = [double]"5.449747"
$ws.Range("I13").Value = [double]"0.09404284435416221"
$ws.Range("J13").Value = [double]"0.0940428443541622"
$ws.Range("M13").Value = [double]"0.08973500000000001"
$ws.Range("N13").Value = [double]"0.269205"
$ws.Range("O13").Value = [double]"0.004645500775423563"
$ws.Range("P13").Value = [double]"0.004645500775423562"
$ws.Range("Q13").Value = [double]"0.1630110156816667"
$ws.Range("R13").Value = [double]"1.467099141135"
$ws.Range("S13").Value = [double]"0.000436876106370298"
$ws.Range("T13").Value = [double]"0.0004368761063702979"

$ws.Range("G14").Value = [double]"0.08973500000000001"
$ws.Range("H14").Value = [double]"0.269205"
$ws.Range("I14").Value = [double]"0.004645500775423563"
$ws.Range("J14").Value = [double]"0.004645500775423562"
$ws.Range("M14").Value = [double]"9.426699666666666"
$ws.Range("N14").Value = [double]"28.280099"
$ws.Range("O14").Value = [double]"0.4880118193702016"
$ws.Range("P14").Value = [double]"0.4880118193702015"
$ws.Range("Q14").Value = [double]"0.8459048945883334"
$ws.Range("R14").Value = [double]"7.613144051295"
$ws.Range("S14").Value = [double]"0.002267059285300136"
$ws.Range("T14").Value = [double]"0.002267059285300135"

$ws.Range("G15").Value = [double]"0.08973500000000001"
$ws.Range("H15").Value = [double]"0.269205"
$ws.Range("I15").Value = [double]"0.004645500775423563"
$ws.Range("J15").Value = [double]"0.004645500775423562"
$ws.Range("M15").Value = [double]"7.983522666666666"
$ws.Range("O15").Value = [double]"0.4132998355002127"
$ws.Range("P15").Value = [double]"0.4132998355002127"
$ws.Range("Q15").Value = [double]"0.7164014064933333"
$ws.Range("R15").Value = [double]"6.44761265844"
$ws.Range("S15").Value = [double]"0.001919984706298669"
$ws.Range("T15").Value = [double]"0.001919984706298669"

$ws.Range("G16").Value = [double]"0.08973500000000001"
$ws.Range("H16").Value = [double]"0.269205"
$ws.Range("I16").Value = [double]"0.004645500775423563"
$ws.Range("J16").Value = [double]"0.004645500775423562"
$ws.Range("M16").Value = [double]"1.816582333333334"
$ws.Range("N16").Value = [double]"5.449747"
$ws.Range("O16").Value = [double]"0.09404284435416221"
$ws.Range("P16").Value = [double]"0.0940428443541622"
$ws.Range("Q16").Value = [double]"0.1630110156816667"
$ws.Range("R16").Value = [double]"1.467099141135"
$ws.Range("S16").Value = [double]"0.000436876106370298"
$ws.Range("T16").Value = [double]"0.0004368761063702979"

$ws.Range("G17").Value = [double]"0.08973500000000001"
$ws.Range("H17").Value = [double]"0.269205"
$ws.Range("I17").Value = [double]"0.004645500775423563"
$ws.Range("J17").Value = [double]"0.004645500775423562"
$ws.Range("M17").Value = [double]"0.08973500000000001"
$ws.Range("N17").Value = [double]"0.269205"
$ws.Range("O17").Value = [double]"0.004645500775423563"
$ws.Range("P17").Value = [double]"0.004645500775423562"
$ws.Range("Q17").Value = [double]"0.008052370225000001"
$ws.Range("R17").Value = [double]"0.07247133202500002"
$ws.Range("S17").Value = [double]"2.158067745446093E-05"
$ws.Range("T17").Value = [double]"2.158067745446092E-05"
